$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (B17): new task entry for the second TensorFlow experiment (image
# classification against the MNIST dataset).
$ws.Range("B17").Value = "4.ทดลองแยกภาพด้วย TensorFlow ใช้ dataSet : MNIST"
$ws.Range("B17").NumberFormat = "@"

# Row 16 (B16): repurpose the old free-form note into a plain "Path : ..." reference
# for the first TensorFlow script.
$ws.Range("B16").Value = "Path : Research/lab/DeepLearning1.py"

# Row 18 (B18): path reference for the new classification script.
$ws.Range("B18").Value = "Path : Research/lab/basic_classification_tShirtSneaker.py"
$ws.Range("B18").NumberFormat = "@"

# Move the active selection down to the new last row, matching where the
# author finished editing.
$ws.Range("B19").Select()
